$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Random, Chronological)
$ws.Range("E3").Value = 9754
$ws.Range("F3").Value = 10794

# Row 4 (Random, Chronological_Opt)
$ws.Range("C4").Value = 0.69
$ws.Range("D4").Value = 0.67
$ws.Range("E4").Value = 13663.7
$ws.Range("F4").Value = 14000

# Row 5 (Random, Random_Opt)
$ws.Range("C5").Value = 0.77
$ws.Range("D5").Value = 0.75
$ws.Range("E5").Value = 10206
$ws.Range("F5").Value = 10689

# Row 6 (GCN, Random)
$ws.Range("C6").Value = 0.76
$ws.Range("E6").Value = 10608
$ws.Range("F6").Value = 11434.8

# Row 7 (GCN, Chronological)
$ws.Range("E7").Value = 7121
$ws.Range("F7").Value = 7741

# Row 8 (GCN, Chronological_Opt)
$ws.Range("E8").Value = 11439
$ws.Range("F8").Value = 11312

# Row 9 (GCN, Random_Opt)
$ws.Range("C9").Value = 0.81
$ws.Range("E9").Value = 7816
$ws.Range("F9").Value = 7988

# Update selection to F11
$ws.Range("F11").Select()
